# Update the LR-pairs data for Tnfsf13b-Tnfrsf13b with refreshed TPM-based values.
# Existing rows 2-5 get new sending/target cluster pairings and recalculated metrics,
# and four additional rows (6-9) are appended to cover all cluster combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Tnfsf13b", "Tnfrsf13b", "ECs", 1, 0.3333333333333333, 0.1273763333333333, 0.382129, 0.02221572899588398, 0.02221572899588398, 1, 0.3333333333333333, 0.4025163333333333, 1.207549, 0.04753623304791887, 0.04753623304791887, 0.05127105464677778, 0.461439491821, 0.001056052070877749, 0.001056052070877749),
    @("ECs", "Tnfsf13b", "Tnfrsf13b", "Resolving-Mac", 1, 0.3333333333333333, 0.1273763333333333, 0.382129, 0.02221572899588398, 0.02221572899588398, 3, 1, 8.065052666666666, 24.195158, 0.9524637669520812, 0.9524637669520811, 1.027296836820222, 9.245671531382, 0.02115967692500623, 0.02115967692500623),
    @("FAPs", "Tnfsf13b", "Tnfrsf13b", "ECs", 3, 1, 4.119982333333334, 12.359947, 0.7185668529619296, 0.7185668529619297, 1, 0.3333333333333333, 0.4025163333333333, 1.207549, 0.04753623304791887, 0.04753623304791887, 1.658360182211444, 14.925241639903, 0.03415796138290794, 0.03415796138290794),
    @("FAPs", "Tnfsf13b", "Tnfrsf13b", "Resolving-Mac", 3, 1, 4.119982333333334, 12.359947, 0.7185668529619296, 0.7185668529619297, 3, 1, 8.065052666666666, 24.195158, 0.9524637669520812, 0.9524637669520811, 33.22787450406955, 299.050870536626, 0.6844088915790217, 0.6844088915790217),
    @("MuSCs", "Tnfsf13b", "Tnfrsf13b", "ECs", 1, 0.3333333333333333, 0.2286643333333333, 0.685993, 0.03988138712600571, 0.03988138712600571, 1, 0.3333333333333333, 0.4025163333333333, 1.207549, 0.04753623304791887, 0.04753623304791887, 0.09204112901744443, 0.8283701611569999, 0.001895810912696079, 0.001895810912696079),
    @("MuSCs", "Tnfsf13b", "Tnfrsf13b", "Resolving-Mac", 1, 0.3333333333333333, 0.2286643333333333, 0.685993, 0.03988138712600571, 0.03988138712600571, 3, 1, 8.065052666666666, 24.195158, 0.9524637669520812, 0.9524637669520811, 1.844189891321556, 16.597709021894, 0.03798557621330963, 0.03798557621330963),
    @("Resolving-Mac", "Tnfsf13b", "Tnfrsf13b", "ECs", 3, 1, 1.257587333333333, 3.772762, 0.2193360309161807, 0.2193360309161807, 1, 0.3333333333333333, 0.4025163333333333, 1.207549, 0.04753623304791887, 0.04753623304791887, 0.5061994422597778, 4.555794980338, 0.0104264086814371, 0.0104264086814371),
    @("Resolving-Mac", "Tnfsf13b", "Tnfrsf13b", "Resolving-Mac", 3, 1, 1.257587333333333, 3.772762, 0.2193360309161807, 0.2193360309161807, 3, 1, 8.065052666666666, 24.195158, 0.9524637669520812, 0.9524637669520811, 10.14250807626622, 91.282572686396, 0.2089096222347437, 0.2089096222347436)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $values[$c - 1]
    }
}
